$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Cell values to update, derived from the diff: a row-by-row price/volume refresh,
# plus two pair swaps: row13<->row14 (Polygon/TRON) and row35<->row36 (Stellar/RenderToken).
$updates = [ordered]@{
    'D2' = '45.642.47'
    'E2' = '  +6.59%  '
    'D3' = '2.385.65'
    'E3' = '  +4.35%  '
    'E4' = '  -0.10%  '
    'D5' = '114.71'
    'E5' = '  +11.37%  '
    'D6' = '317.81'
    'E6' = '  +2.47%  '
    'E7' = '  +2.05%  '
    'E8' = '  -0.23%  '
    'D9' = '0.630'
    'E9' = '  +4.79%  '
    'D10' = '43.07'
    'E10' = '  +11.55%  '
    'D11' = '0.0939'
    'E11' = '  +4.42%  '
    'D12' = '8.72'
    'E12' = '  +6.45%  '
    'B13' = 'TRON'
    'C13' = 'https://coinranking.com/coin/qUhEFk1I61atv+tron-trx'
    'D13' = '0.109'
    'E13' = '  +1.14%  '
    'B14' = 'Polygon'
    'C14' = 'https://coinranking.com/coin/uW2tk-ILY0ii+polygon-matic'
    'D14' = '1.01'
    'E14' = '  +4.80%  '
    'D15' = '15.94'
    'E15' = '  +4.75%  '
    'D16' = '2.745.27'
    'E16' = '  +4.29%  '
    'D17' = '2.378.85'
    'E17' = '  +4.13%  '
    'D18' = '45.506.45'
    'E18' = '  +6.66%  '
    'D19' = '7.58'
    'E19' = '  +4.25%  '
    'E20' = '  +3.93%  '
    'D21' = '13.35'
    'E21' = '  -0.16%  '
    'D22' = '74.97'
    'E22' = '  +2.52%  '
    'E23' = '  +4.16%  '
    'D24' = '269.13'
    'E24' = '  +0.23%  '
    'E25' = '  +9.43%  '
    'E26' = '  -0.62%  '
    'D27' = '7.67'
    'E27' = '  +11.10%  '
    'D28' = '11.33'
    'E28' = '  +5.75%  '
    'E29' = '  +1.34%  '
    'D30' = '22.93'
    'E30' = '  +2.92%  '
    'D31' = '38.92'
    'E31' = '  +8.82%  '
    'D32' = '0.0969'
    'E32' = '  +14.79%  '
    'D33' = '171.16'
    'E33' = '  +4.14%  '
    'D34' = '2.98'
    'E34' = '  +17.41%  '
    'B35' = 'RenderToken'
    'C35' = 'https://coinranking.com/coin/7C4Mh4xy1yDel+rendertoken-rndr'
    'D35' = '5.00'
    'E35' = '  +11.16%  '
    'B36' = 'Stellar'
    'C36' = 'https://coinranking.com/coin/f3iaFeCKEmkaZ+stellar-xlm'
    'D36' = '0.132'
    'E36' = '  +1.84%  '
    'D37' = '0.119'
    'E37' = '  +7.64%  '
    'D38' = '4.12'
    'E38' = '  +14.49%  '
    'D39' = '3.06'
    'E39' = '  +11.43%  '
    'D40' = '0.0367'
    'E40' = '  +6.42%  '
    'D41' = '1.73'
    'E41' = '  +11.14%  '
    'D42' = '104.19'
    'E42' = '  -6.99%  '
    'D43' = '0.240'
    'E43' = '  +7.03%  '
    'D44' = '71.49'
    'E44' = '  +2.45%  '
    'D45' = '13.32'
    'E45' = '  +10.84%  '
    'E46' = '  -0.55%  '
    'E47' = '  +12.11%  '
    'D48' = '116.32'
    'E48' = '  +5.93%  '
    'E49' = '  +17.01%  '
    'D50' = '9.37'
    'E50' = '  +8.35%  '
    'D51' = '79.41'
    'E51' = '  +3.19%  '
}

foreach ($addr in $updates.Keys) {
    $value = $updates[$addr]
    $range = $ws.Range($addr)
    if ($value -match '^-?\d+(\.\d+)?$') {
        # The new value looks like a plain number (e.g. "114.71"). The sheet stores
        # these as text, so force text storage - otherwise Excel would silently
        # convert the string into a float and lose exact formatting/trailing zeros.
        $range.NumberFormat = "@"
        $range.Value = $value
        $range.Style = "Normal"
    } else {
        $range.Value = $value
    }
}
